$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 28 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 924.46155
$ws.Range("I11").Value = 924.46155
$ws.Range("K11").Value = 924.46155
$ws.Range("M11").Value = -784.46155
$ws.Range("H51").Value = 32667.666
$ws.Range("J51").Value = 11277.25
$ws.Range("L51").Value = 11277.25
$ws.Range("N51").Value = -12245.25
$ws.Range("H80").Value = 759.8
$ws.Range("I80").Value = 699.6667
$ws.Range("K80").Value = 2099.0001
$ws.Range("M80").Value = -1101.0001
$ws.Range("H83").Value = 759.8
$ws.Range("I83").Value = 699.6667
$ws.Range("K83").Value = 6297.0003
$ws.Range("M83").Value = -1305.0003
$ws.Range("H107").Value = 2795.125
$ws.Range("I107").Value = 2286.3845
$ws.Range("K107").Value = 2286.3845
$ws.Range("M107").Value = -366.3845000000001
$ws.Range("H137").Value = 2160.0557
$ws.Range("J137").Value = 2235.2144
$ws.Range("L137").Value = 6705.6432
$ws.Range("N137").Value = -11805.6432
$ws.Range("H141").Value = 9636.333000000001
$ws.Range("I141").Value = 6035.875
$ws.Range("K141").Value = 18107.625
$ws.Range("M141").Value = -12927.625

# --- Sheet ARM: 27 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1632
$ws.Range("I2").Value = 1670.7142
$ws.Range("K2").Value = 1670.7142
$ws.Range("M2").Value = -1557.7142
$ws.Range("H8").Value = 344333
$ws.Range("J8").Value = 29999
$ws.Range("L8").Value = 29999
$ws.Range("N8").Value = -30287
$ws.Range("H32").Value = 2042.9584
$ws.Range("I32").Value = 1881.8695
$ws.Range("J32").Value = 5748
$ws.Range("K32").Value = 1881.8695
$ws.Range("L32").Value = 5748
$ws.Range("M32").Value = -1594.8695
$ws.Range("N32").Value = -6322
$ws.Range("H63").Value = 1165.3334
$ws.Range("I63").Value = 1223.625
$ws.Range("K63").Value = 1223.625
$ws.Range("M63").Value = -537.625
$ws.Range("H66").Value = 1165.3334
$ws.Range("I66").Value = 1223.625
$ws.Range("K66").Value = 6118.125
$ws.Range("M66").Value = -2686.125
$ws.Range("H116").Value = 1632
$ws.Range("I116").Value = 1670.7142
$ws.Range("K116").Value = 1670.7142
$ws.Range("M116").Value = 623.2858000000001

# --- Sheet BSM: 20 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1632
$ws.Range("I3").Value = 1670.7142
$ws.Range("K3").Value = 1670.7142
$ws.Range("M3").Value = -1556.7142
$ws.Range("H99").Value = 3303.9412
$ws.Range("I99").Value = 3277.8667
$ws.Range("K99").Value = 3277.8667
$ws.Range("M99").Value = -1779.8667
$ws.Range("H105").Value = 2873.9333
$ws.Range("I105").Value = 2716
$ws.Range("K105").Value = 2716
$ws.Range("M105").Value = -969
$ws.Range("H132").Value = 65249.5
$ws.Range("J132").Value = 65249.5
$ws.Range("L132").Value = 65249.5
$ws.Range("N132").Value = -75369.5
$ws.Range("H134").Value = 15218.2
$ws.Range("I134").Value = 27425
$ws.Range("K134").Value = 82275
$ws.Range("M134").Value = -79740

# --- Sheet CRP: 38 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 949.6667
$ws.Range("I22").Value = 449.5
$ws.Range("K22").Value = 449.5
$ws.Range("M22").Value = -99.5
$ws.Range("H31").Value = 2635
$ws.Range("I31").Value = 1481.6428
$ws.Range("K31").Value = 1481.6428
$ws.Range("M31").Value = -1186.6428
$ws.Range("H34").Value = 2635
$ws.Range("I34").Value = 1481.6428
$ws.Range("K34").Value = 1481.6428
$ws.Range("M34").Value = -1279.6428
$ws.Range("H99").Value = 3338.5
$ws.Range("I99").Value = 2925
$ws.Range("K99").Value = 2925
$ws.Range("M99").Value = -1427
$ws.Range("H105").Value = 1326.1428
$ws.Range("I105").Value = 1255.4286
$ws.Range("J105").Value = 1467.5714
$ws.Range("K105").Value = 1255.4286
$ws.Range("L105").Value = 1467.5714
$ws.Range("M105").Value = 491.5714
$ws.Range("N105").Value = -4961.5714
$ws.Range("H122").Value = 3825.3
$ws.Range("I122").Value = 3530.875
$ws.Range("K122").Value = 10592.625
$ws.Range("M122").Value = -8142.625
$ws.Range("H126").Value = 3338.5
$ws.Range("I126").Value = 2925
$ws.Range("K126").Value = 8775
$ws.Range("M126").Value = -6305
$ws.Range("H134").Value = 5942.593
$ws.Range("I134").Value = 5190.5884
$ws.Range("J134").Value = 7221
$ws.Range("K134").Value = 15571.7652
$ws.Range("L134").Value = 21663
$ws.Range("M134").Value = -13036.7652
$ws.Range("N134").Value = -26733

# --- Sheet CUL: 34 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1045409.7
$ws.Range("I11").Value = 1567964.5
$ws.Range("K11").Value = 4703893.5
$ws.Range("M11").Value = -4703753.5
$ws.Range("H22").Value = 6184.125
$ws.Range("I22").Value = 4912.1665
$ws.Range("K22").Value = 14736.4995
$ws.Range("M22").Value = -14567.4995
$ws.Range("H27").Value = 6184.125
$ws.Range("I27").Value = 4912.1665
$ws.Range("K27").Value = 14736.4995
$ws.Range("M27").Value = -14634.4995
$ws.Range("H59").Value = 2499
$ws.Range("I59").Value = 2499
$ws.Range("K59").Value = 7497
$ws.Range("M59").Value = -6957
$ws.Range("H60").Value = 1338.4
$ws.Range("I60").Value = 1338.4
$ws.Range("K60").Value = 4015.2
$ws.Range("M60").Value = -3764.2
$ws.Range("H86").Value = 4243
$ws.Range("I86").Value = 371.2
$ws.Range("J86").Value = 8114.8
$ws.Range("K86").Value = 1113.6
$ws.Range("L86").Value = 24344.4
$ws.Range("M86").Value = 72.40000000000009
$ws.Range("N86").Value = -26716.4
$ws.Range("H89").Value = 4243
$ws.Range("I89").Value = 371.2
$ws.Range("J89").Value = 8114.8
$ws.Range("K89").Value = 3340.8
$ws.Range("L89").Value = 73033.2
$ws.Range("M89").Value = 2587.2
$ws.Range("N89").Value = -84889.2

# --- Sheet GSM: 23 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4651.5835
$ws.Range("I102").Value = 4398.6875
$ws.Range("K102").Value = 4398.6875
$ws.Range("M102").Value = -2776.6875
$ws.Range("H122").Value = 4599.1665
$ws.Range("I122").Value = 3866
$ws.Range("K122").Value = 11598
$ws.Range("M122").Value = -9148
$ws.Range("H126").Value = 5191.6924
$ws.Range("I126").Value = 2747.5
$ws.Range("K126").Value = 8242.5
$ws.Range("M126").Value = -5772.5
$ws.Range("H132").Value = 2357.7144
$ws.Range("I132").Value = 2401.3684
$ws.Range("J132").Value = 1943
$ws.Range("K132").Value = 7204.1052
$ws.Range("L132").Value = 5829
$ws.Range("M132").Value = -4674.1052
$ws.Range("N132").Value = -10889
$ws.Range("H135").Value = 113000
$ws.Range("J135").Value = 113000
$ws.Range("L135").Value = 113000
$ws.Range("N135").Value = -123140

# --- Sheet LTW: 30 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1250.1428
$ws.Range("I40").Value = 1222.65
$ws.Range("K40").Value = 1222.65
$ws.Range("M40").Value = -1086.65
$ws.Range("H46").Value = 1585.6111
$ws.Range("I46").Value = 928.3333
$ws.Range("K46").Value = 928.3333
$ws.Range("M46").Value = -740.3333
$ws.Range("H55").Value = 606.3333
$ws.Range("I55").Value = 300
$ws.Range("J55").Value = 667.6
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 667.6
$ws.Range("M55").Value = -127
$ws.Range("N55").Value = -1013.6
$ws.Range("H76").Value = 38999.668
$ws.Range("J76").Value = 38999.668
$ws.Range("L76").Value = 38999.668
$ws.Range("N76").Value = -39675.668
$ws.Range("H79").Value = 38999.668
$ws.Range("J79").Value = 38999.668
$ws.Range("L79").Value = 38999.668
$ws.Range("N79").Value = -41339.668
$ws.Range("H132").Value = 5352.231
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 5689
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 17067
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -22127

# --- Sheet WVR: 23 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 15200
$ws.Range("I37").Value = 15200
$ws.Range("K37").Value = 15200
$ws.Range("M37").Value = -14997
$ws.Range("H113").Value = 648.25
$ws.Range("I113").Value = 465
$ws.Range("K113").Value = 1395
$ws.Range("M113").Value = 775
$ws.Range("H122").Value = 4220.6206
$ws.Range("I122").Value = 2662.7334
$ws.Range("K122").Value = 7988.2002
$ws.Range("M122").Value = -5538.2002
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3250
$ws.Range("K126").Value = 9750
$ws.Range("M126").Value = -7280
$ws.Range("H132").Value = 2347.7896
$ws.Range("I132").Value = 2054.3125
$ws.Range("J132").Value = 3913
$ws.Range("K132").Value = 6162.9375
$ws.Range("L132").Value = 11739
$ws.Range("M132").Value = -3632.9375
$ws.Range("N132").Value = -16799

Write-Output "Applied 223 cell updates across 8 sheets"